# Estabilizacion pago de tarjetas de credito
# Updates the data-driven test rows for TarjetaCreditoPropia:
#  - usuario (D2:D4):      userrobot1     -> autotest10
#  - clave   (E2:E4):      6789           -> 1234
#  - tipoCuenta (T2):      Ahorros        -> Corriente
#  - numeroCuenta (U2):    406-739440-03  -> 406-125170-00
#  - numeroCuenta (U3:U4): 406-739440-03  -> 406-725170-06

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (usuario): userrobot1 -> autotest10.
# Leading apostrophe + explicit font name reproduces the "quote prefix"
# text style Excel applies to this column (numeric-looking text value).
$ws.Range("D2").Value = "'autotest10"
$ws.Range("D2").Font.Name = "Calibri"
$ws.Range("D3").Value = "'autotest10"
$ws.Range("D3").Font.Name = "Calibri"
$ws.Range("D4").Value = "'autotest10"
$ws.Range("D4").Font.Name = "Calibri"

# Column E (clave): 6789 -> 1234 (keep its existing quote-prefix text style)
$ws.Range("E2").Value2 = "'1234"
$ws.Range("E3").Value2 = "'1234"
$ws.Range("E4").Value2 = "'1234"

# Column T row 2 only (tipoCuenta): Ahorros -> Corriente
$ws.Range("T2").Value2 = "Corriente"

# Column U (numeroCuenta)
$ws.Range("U2").Value2 = "406-125170-00"
$ws.Range("U3").Value2 = "406-725170-06"
$ws.Range("U4").Value2 = "406-725170-06"

# Sheet view changes: scroll so column K is the left-most visible column,
# and select I2 (mirrors topLeftCell="K1" / selection activeCell="I2")
$excel.ActiveWindow.ScrollColumn = 11
$excel.ActiveWindow.ScrollRow = 1
$null = $ws.Range("I2").Select()
